# Add the two new tweets as new shared-string rows at the bottom of column A
# (rows 35 and 36), matching the new uniqueCount=36 sharedStrings.xml entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow35 = "Öğrencinin online derste kamera açması zorunlu bile değilken zorla kamera açtırıp,online eğitim sürecinde yoklama alınmayacağını bildirip,yoklamayla tehdit etmenizin üstüne bide uzaktan gözetmeli sınav mı çıkartiyorsunuz? Işık uyuma! #isikbizimledegil"
$newRow36 = "Günlerdir sesimizi üniversitemize duyuramıyoruz. Online dersler ve sınavlarla ilgili adaletsizliğeve düşüncesizliğe karşı duruyoruz fakat kaideye alınmıyoruz.Lütfen sesimiz olun. #soruisaretleriyle #isikunısesver @fatihportakal"

$ws.Range("A35").Value = $newRow35
$ws.Range("A36").Value = $newRow36

# Mirror the refreshed view state from the saved workbook: the window had
# scrolled down a bit and the selection moved past the new last row.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B39").Select()
